$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -11.13889999999999
$ws.Range("A3").Value = -21.45610000000002
$ws.Range("C5").Value = -14.30419999999999
$ws.Range("D5").Value = -8.261699999999992
$ws.Range("E7").Value = 11.9466
$ws.Range("D9").Value = -8.738700000000001
$ws.Range("D11").Value = -8.727099999999998
$ws.Range("E11").Value = 12.61369999999998
$ws.Range("A14").Value = -20.53149999999998
$ws.Range("E19").Value = 13.3972
$ws.Range("A21").Value = -21.51520000000002
$ws.Range("D21").Value = -8.026100000000008
$ws.Range("E21").Value = 13.24560000000002
$ws.Range("A23").Value = -21.37330000000003
$ws.Range("A25").Value = -22.62070000000004
